# Apply cryptos.xlsx price/volume update (commit: "Updated cryptos list on Sat Jul 29 03:27:26 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.384.37"
$ws.Range("E2").Value = "  +0.56%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.874.60"
$ws.Range("E3").Value = "  +0.68%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7120"
$ws.Range("E5").Value = "  -0.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.08"
$ws.Range("E6").Value = "  +0.62%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3112"
$ws.Range("E8").Value = "  +0.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07787"
$ws.Range("E9").Value = "  +1.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.14"
$ws.Range("E10").Value = "  +0.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08441"
$ws.Range("E11").Value = "  +1.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.865.58"
$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.241"
$ws.Range("E13").Value = "  +0.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7124"
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.25"
$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.386.13"
$ws.Range("E16").Value = "  +0.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008239"
$ws.Range("E17").Value = "  +5.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.042"
$ws.Range("E18").Value = "  +0.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.21"
$ws.Range("E19").Value = "  -1.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.27"
$ws.Range("E20").Value = "  +0.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.120.39"
$ws.Range("E21").Value = "  -0.94%  "

$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.787"
$ws.Range("E23").Value = "  -2.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9999"
$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1607"
$ws.Range("E25").Value = "  -0.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.62"
$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.064"
$ws.Range("E27").Value = "  +1.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.47"
$ws.Range("E28").Value = "  -0.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.511"
$ws.Range("E29").Value = "  +0.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.428"
$ws.Range("E30").Value = "  -0.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.288"
$ws.Range("E31").Value = "  -4.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.305"
$ws.Range("E32").Value = "  +1.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05286"
$ws.Range("E33").Value = "  +1.98%  "

$ws.Range("E34").Value = "  +0.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.179"
$ws.Range("E35").Value = "  +0.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7445"
$ws.Range("E36").Value = "  -8.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.696"
$ws.Range("E37").Value = "  +0.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01873"
$ws.Range("E38").Value = "  +0.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.217.06"
$ws.Range("E39").Value = "  +4.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.725"
$ws.Range("E40").Value = "  +1.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.474"
$ws.Range("E41").Value = "  +4.01%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8877"
$ws.Range("E42").Value = "  -1.92%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.82"
$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.46"
$ws.Range("E44").Value = "  +7.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.019.38"
$ws.Range("E46").Value = "  -1.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.816"
$ws.Range("E47").Value = "  +1.89%  "

$ws.Range("E48").Value = "  +0.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000123"
$ws.Range("E49").Value = "  +1.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.376"
$ws.Range("E50").Value = "  -0.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4323"
$ws.Range("E51").Value = "  +0.92%  "
